# Add a new "Complaint" participant-assignment rule row to Sheet1, mirroring
# the existing "Case File" row (row 20) but targeting the COMPLAINT object
# type instead of CASE_FILE.
#
# Commit message: "File changed so that now we have mechanism to prevent
# combination for Assignee and No Access participant for a Complaint"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make sure we're working on the right sheet / it stays the active one.
$ws.Activate()

# Copy the formatting (fills/borders/number formats) of the existing
# "Case File" rule row (row 20) down onto the new row 21 so the new row
# visually matches the rest of the rule table.
$ws.Range("B20:F20").Copy()
$ws.Range("B21:F21").PasteSpecial(-4122)

# Fill in the new rule's values. Column C differentiates it from row 20
# (COMPLAINT instead of CASE_FILE); the rest of the logic (columns B, D, E,
# F) mirrors the Case File rule, just re-worded for "Complaint".
$ws.Range("B21").Value = "Complaint - Check participants list for NoAccess & Owner"
$ws.Range("C21").Value = "COMPLAINT"
$ws.Range("D21").Value = $ws.Range("D20").Value()
$ws.Range("E21").Value = $ws.Range("E20").Value()
$ws.Range("F21").Value = $ws.Range("F20").Value()

# Match row 20's row height (wrapped rule-name text needs the extra height).
$ws.Rows.Item(21).RowHeight = 45

# Update the view so the active cell follows the new last row, same as the
# author's selection after adding the row.
$ws.Range("D21").Select() | Out-Null
